$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.908.41"
$ws.Range("E2").Value = "  +0.87%  "

$ws.Range("D3").Value = "2.347.99"
$ws.Range("E3").Value = "  +1.25%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'539.33"
$ws.Range("E5").Value = "  +1.80%  "

$ws.Range("D6").Value = "'135.65"
$ws.Range("E6").Value = "  +2.46%  "

$ws.Range("E7").Value = "  +0.33%  "

$ws.Range("E9").Value = "  +0.48%  "

$ws.Range("E10").Value = "  +4.59%  "

$ws.Range("E11").Value = "  -0.71%  "

$ws.Range("E12").Value = "  +1.88%  "

$ws.Range("D13").Value = "'23.78"
$ws.Range("E13").Value = "  +1.54%  "

$ws.Range("D14").Value = "2.766.79"
$ws.Range("E14").Value = "  +1.06%  "

$ws.Range("D15").Value = "57.901.01"
$ws.Range("E15").Value = "  +1.14%  "

$ws.Range("E16").Value = "  +0.57%  "

$ws.Range("D17").Value = "2.319.68"
$ws.Range("E17").Value = "  -0.75%  "

$ws.Range("E18").Value = "  +2.57%  "

$ws.Range("D19").Value = "'332.30"
$ws.Range("E19").Value = "  -1.20%  "

$ws.Range("E20").Value = "  +2.79%  "

$ws.Range("E21").Value = "  -0.85%  "

$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.12%  "

$ws.Range("D23").Value = "'62.88"
$ws.Range("E23").Value = "  +1.99%  "

$ws.Range("E24").Value = "  -0.52%  "

$ws.Range("E25").Value = "  -2.36%  "

$ws.Range("E26").Value = "  +0.57%  "

$ws.Range("E27").Value = "  +0.79%  "

$ws.Range("D28").Value = "'171.97"

$ws.Range("E29").Value = "  +1.53%  "

$ws.Range("E30").Value = "  +1.64%  "

$ws.Range("D31").Value = "'6.14"
$ws.Range("E31").Value = "  +0.48%  "

$ws.Range("E32").Value = "  +10.50%  "

$ws.Range("D33").Value = "'18.45"
$ws.Range("E33").Value = "  -0.26%  "

$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("E37").Value = "  -0.40%  "

$ws.Range("D38").Value = "'1.64"
$ws.Range("E38").Value = "  +4.09%  "

$ws.Range("D39").Value = "'39.25"
$ws.Range("E39").Value = "  +0.24%  "

$ws.Range("D40").Value = "'145.74"
$ws.Range("E40").Value = "  -2.23%  "

$ws.Range("D41").Value = "'294.50"
$ws.Range("E41").Value = "  +4.52%  "

$ws.Range("E42").Value = "  +0.70%  "

$ws.Range("E43").Value = "  +1.24%  "

$ws.Range("E44").Value = "  +1.85%  "

$ws.Range("D45").Value = "'19.28"
$ws.Range("E45").Value = "  +2.24%  "

$ws.Range("E46").Value = "  +0.43%  "

$ws.Range("E47").Value = "  +0.79%  "

$ws.Range("E48").Value = "  +1.65%  "

$ws.Range("D49").Value = "'0.385"
$ws.Range("E49").Value = "  +0.90%  "

$ws.Range("D50").Value = "'17.47"
$ws.Range("E50").Value = "  -0.32%  "

$ws.Range("D51").Value = "'11.08"
$ws.Range("E51").Value = "  +0.51%  "

# Row 35/36 swap: FirstDigitalUSD <-> NEARProtocol with updated values
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "'4.22"
$ws.Range("E35").Value = "  +5.87%  "

$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  +0.22%  "
